# feat: add multi-institute support
#
# The "Siswa" (student) import template reorders its columns: no_induk,
# jenis_kelamin and nisn move up right after "nama" (and before
# tempat_lahir/tgl_lahir), and the address block (alamat, rt, rw,
# kelurahan, kecamatan, kab_kota, provinsi) moves up to sit right after
# tgl_lahir instead of after the parents' block.
#
# We reorder the existing columns (preserving per-column width/content)
# using Cut + Insert so the column formatting (width) travels with the
# data instead of being re-typed with a lossy numeric width.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move column D (no_induk) to B.
$ws.Columns.Item(4).Cut() | Out-Null
$ws.Columns.Item(2).Insert() | Out-Null

# Move column F (jenis_kelamin) to C.
$ws.Columns.Item(6).Cut() | Out-Null
$ws.Columns.Item(3).Insert() | Out-Null

# Move column F (nisn, shifted after previous moves) to D.
$ws.Columns.Item(6).Cut() | Out-Null
$ws.Columns.Item(4).Insert() | Out-Null

# Move column K (alamat) to G.
$ws.Columns.Item(11).Cut() | Out-Null
$ws.Columns.Item(7).Insert() | Out-Null

# Move column L (rt) to H.
$ws.Columns.Item(12).Cut() | Out-Null
$ws.Columns.Item(8).Insert() | Out-Null

# Move column M (rw) to I.
$ws.Columns.Item(13).Cut() | Out-Null
$ws.Columns.Item(9).Insert() | Out-Null

# Move column N (kelurahan) to J.
$ws.Columns.Item(14).Cut() | Out-Null
$ws.Columns.Item(10).Insert() | Out-Null

# Move column O (kecamatan) to K.
$ws.Columns.Item(15).Cut() | Out-Null
$ws.Columns.Item(11).Insert() | Out-Null

# Move column P (kab_kota) to L.
$ws.Columns.Item(16).Cut() | Out-Null
$ws.Columns.Item(12).Insert() | Out-Null

# Move column Q (provinsi) to M.
$ws.Columns.Item(17).Cut() | Out-Null
$ws.Columns.Item(13).Insert() | Out-Null

# Final column order is now:
# A nama | B no_induk | C jenis_kelamin | D nisn | E tempat_lahir |
# F tgl_lahir | G alamat | H rt | I rw | J kelurahan | K kecamatan |
# L kab_kota | M provinsi | N nama_ayah | O pekerjaan_ayah |
# P nama_ibu | Q pekerjaan_ibu

# Update the saved selection to match the authored workbook (O16).
$ws.Range("O16").Select() | Out-Null
